$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of the column to the left (M) before inserting, so the
# newly inserted column inherits the same display width.
$leftWidth = $ws.Columns("M").ColumnWidth

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Make this the active sheet/tab, with the given cell selected.
$ws.Activate()
$ws.Range("M14").Select() | Out-Null
